$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) data rows 2-19 keyed by player name (column A),
# so the row order can be rearranged without losing any associated data.
$rows = @{}
for ($r = 2; $r -le 19; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    $pos  = $ws.Cells.Item($r, 2).Value()
    $team = $ws.Cells.Item($r, 3).Value()
    $rows[$name] = @($name, $pos, $team)
}

# New desired order of players (rows 2-19)
$newOrder = @(
    "Donovan Mitchell",
    "Bradley Beal",
    "Michael Porter Jr.",
    "Deandre Ayton",
    "Myles Turner",
    "Victor Wembanyama",
    "Malik Beasley",
    "Payton Pritchard",
    "Kristaps Porzingis",
    "Domantas Sabonis",
    "Tari Eason",
    "De'Andre Hunter",
    "Josh Hart",
    "Dyson Daniels",
    "Jamal Murray",
    "Robert Williams III",
    "Cam Thomas",
    "P.J. Washington"
)

$r = 2
foreach ($name in $newOrder) {
    $data = $rows[$name]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $r = $r + 1
}
